# añadiento archivos de prueba manual
#
# Corrects the wording in a few cells of the manual-test-case sheet:
#  - D7:  "resultado esperado" text for the results-count test case
#  - A11: adds the missing space in the search phrase for the
#         "no results" negative test case
#  - A12: fixes the "no e muestren" -> "no se muestren" typo
#  - D12: rewrites the "no results found" system message

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "la catidad de resultados se muestran de manera correcta junto con el tiempo que se tardó el sistema en realizar la búsqueda"

$ws.Range("A11").Value = "Ingresar la palabra asccffaxadd agdccssessssaefjkk,l en la barra de búsqueda y presionar la tecla enter"

$ws.Range("A12").Value = "Validar que no se muestren resultados  de  búsqueda"

$ws.Range("D12").Value = "El sistema debe mostrar un mensaje que no se han encontrado datos.
Sugerencias:"

# Leave the selection where the author left it before saving.
$null = $ws.Range("I12").Select()
